$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4293483371413345
$ws.Range("C2").Value = 0.1877078620774881
$ws.Range("D2").Value = 0.04969493119833857
$ws.Range("F2").Value = 0.9677544571148289
$ws.Range("G2").Value = 0.002465892824517755
$ws.Range("I2").Value = 0.9550838363352696
$ws.Range("K2").Value = 0.2626602799647628
$ws.Range("L2").Value = 0.304282200983792
$ws.Range("M2").Value = 0.1644841351855959
$ws.Range("O2").Value = 3.466989768885725
# Row 3
$ws.Range("B3").Value = 0.3935028817097646
$ws.Range("C3").Value = 0.1867812462712806
$ws.Range("D3").Value = 0.04743121432228747
$ws.Range("F3").Value = 0.969700940399477
$ws.Range("G3").Value = 0.002468236298011073
$ws.Range("I3").Value = 0.9625522134310422
$ws.Range("K3").Value = 0.2296939974483507
$ws.Range("L3").Value = 0.3010290703549146
$ws.Range("M3").Value = 0.1570552011898378
$ws.Range("O3").Value = 3.486902580698526
# Row 4
$ws.Range("B4").Value = 0.3715769101279136
$ws.Range("C4").Value = 0.1862153439656709
$ws.Range("D4").Value = 0.04602624700727631
$ws.Range("F4").Value = 0.9714272709228453
$ws.Range("G4").Value = 0.002469752834022627
$ws.Range("I4").Value = 0.9675783278341079
$ws.Range("K4").Value = 0.2094071574747858
$ws.Range("L4").Value = 0.2991861802454849
$ws.Range("M4").Value = 0.152559819005134
$ws.Range("O4").Value = 3.500883428954609
# Row 5
$ws.Range("B5").Value = 0.3626634378722144
$ws.Range("C5").Value = 0.1859855278985236
$ws.Range("D5").Value = 0.04544995625015247
$ws.Range("F5").Value = 0.9722644031684951
$ws.Range("G5").Value = 0.00247039041328947
$ws.Range("I5").Value = 0.9697373186261444
$ws.Range("K5").Value = 0.2011291582608976
$ws.Range("L5").Value = 0.298474125870797
$ws.Range("M5").Value = 0.1507446402707693
$ws.Range("O5").Value = 3.507022010423128
# Row 6
$ws.Range("B6").Value = 0.3611846811118085
$ws.Range("C6").Value = 0.1859474158920307
$ws.Range("D6").Value = 0.04535403756504763
$ws.Range("F6").Value = 0.9724114817924061
$ws.Range("G6").Value = 0.0024704974671645
$ws.Range("I6").Value = 0.9701025111593715
$ws.Range("K6").Value = 0.1997539555182897
$ws.Range("L6").Value = 0.2983582438176811
$ws.Range("M6").Value = 0.1504442451592993
$ws.Range("O6").Value = 3.508067973056669
# Row 7
$ws.Range("B7").Value = 0.37145661191974
$ws.Range("C7").Value = 0.186212241333763
$ws.Range("D7").Value = 0.04601849011897485
$ws.Range("F7").Value = 0.9714380195936485
$ws.Range("G7").Value = 0.002469761353293476
$ws.Range("I7").Value = 0.9676069960265394
$ws.Range("K7").Value = 0.2092955611815484
$ws.Range("L7").Value = 0.299176419469525
$ws.Range("M7").Value = 0.1525352709790617
$ws.Range("O7").Value = 3.500964429135109
# Row 8
$ws.Range("B8").Value = 0.4169718849896071
$ws.Range("C8").Value = 0.1873877499223937
$ws.Range("D8").Value = 0.04891753716918146
$ws.Range("F8").Value = 0.9683153968956475
$ws.Range("G8").Value = 0.002466684775632662
$ws.Range("I8").Value = 0.9575675313122574
$ws.Range("K8").Value = 0.2513032028869304
$ws.Range("L8").Value = 0.303128489923381
$ws.Range("M8").Value = 0.1619090108709713
$ws.Range("O8").Value = 3.473491681591582
# Row 9
$ws.Range("B9").Value = 0.5068663648028178
$ws.Range("C9").Value = 0.1897159465252827
$ws.Range("D9").Value = 0.05448236403120177
$ws.Range("F9").Value = 0.9664046299458562
$ws.Range("G9").Value = 0.002461265004581262
$ws.Range("I9").Value = 0.9413736944611699
$ws.Range("K9").Value = 0.3333029934676404
$ws.Range("L9").Value = 0.3121021201440897
$ws.Range("M9").Value = 0.1808101756588982
$ws.Range("O9").Value = 3.433532960790302
# Row 10
$ws.Range("B10").Value = 0.5732791535059505
$ws.Range("C10").Value = 0.1914391889294009
$ws.Range("D10").Value = 0.05849672287616414
$ws.Range("F10").Value = 0.9675677560779974
$ws.Range("G10").Value = 0.002457653391130337
$ws.Range("I10").Value = 0.931603822935223
$ws.Range("K10").Value = 0.3933016433873036
$ws.Range("L10").Value = 0.3194386015944843
$ws.Range("M10").Value = 0.1950089492977369
$ws.Range("O10").Value = 3.412655190790844
# Row 11
$ws.Range("B11").Value = 0.6035674961917721
$ws.Range("C11").Value = 0.1922256317076432
$ws.Range("D11").Value = 0.06030670603875876
$ws.Range("F11").Value = 0.9686541249008371
$ws.Range("G11").Value = 0.002456090004917989
$ws.Range("I11").Value = 0.9276208771087795
$ws.Range("K11").Value = 0.4205398541262468
$ws.Range("L11").Value = 0.3229371703184114
$ws.Range("M11").Value = 0.2015352502464296
$ws.Range("O11").Value = 3.404998208666029
# Row 12
$ws.Range("B12").Value = 0.6150474389937415
$ws.Range("C12").Value = 0.1925237753189322
$ws.Range("D12").Value = 0.06098975280357877
$ws.Range("F12").Value = 0.969145595551538
$ws.Range("G12").Value = 0.002455509372598929
$ws.Range("I12").Value = 0.9261789601889845
$ws.Range("K12").Value = 0.4308458796699028
$ws.Range("L12").Value = 0.3242851015850903
$ws.Range("M12").Value = 0.2040161482122897
$ws.Range("O12").Value = 3.40236328592681
# Row 13
$ws.Range("B13").Value = 0.6125745753871286
$ws.Range("C13").Value = 0.1924595503306819
$ws.Range("D13").Value = 0.06084275165189723
$ws.Range("F13").Value = 0.9690361874435069
$ws.Range("G13").Value = 0.002455633916470492
$ws.Range("I13").Value = 0.9264865525958896
$ws.Range("K13").Value = 0.4286266786382953
$ws.Range("L13").Value = 0.3239937743447712
$ws.Range("M13").Value = 0.2034814207706788
$ws.Range("O13").Value = 3.402918995725145
# Row 14
$ws.Range("B14").Value = 0.6045117531498363
$ws.Range("C14").Value = 0.1922501536177847
$ws.Range("D14").Value = 0.06036294806519749
$ws.Range("F14").Value = 0.9686929537075315
$ws.Range("G14").Value = 0.002456042007999544
$ws.Range("I14").Value = 0.9275009202944418
$ws.Range("K14").Value = 0.4213879102300666
$ws.Range("L14").Value = 0.3230476029305862
$ws.Range("M14").Value = 0.2017391650039784
$ws.Range("O14").Value = 3.40477612903797
# Row 15
$ws.Range("B15").Value = 0.5995743776515212
$ws.Range("C15").Value = 0.1921219348252237
$ws.Range("D15").Value = 0.06006874700295839
$ws.Range("F15").Value = 0.9684931413419804
$ws.Range("G15").Value = 0.002456293457303058
$ws.Range("I15").Value = 0.9281308886042439
$ws.Range("K15").Value = 0.4169528341303135
$ws.Range("L15").Value = 0.3224710515487033
$ws.Range("M15").Value = 0.2006732200496444
$ws.Range("O15").Value = 3.405948135664829
# Row 16
$ws.Range("B16").Value = 0.5713012300898583
$ws.Range("C16").Value = 0.1913878414785799
$ws.Range("D16").Value = 0.05837810833818224
$ws.Range("F16").Value = 0.9675079688032042
$ws.Range("G16").Value = 0.002457757158555268
$ws.Range("I16").Value = 0.9318734015997876
$ws.Range("K16").Value = 0.391520401770066
$ws.Range("L16").Value = 0.3192131984659312
$ws.Range("M16").Value = 0.1945837810361084
$ws.Range("O16").Value = 3.413192623155794
# Row 17
$ws.Range("B17").Value = 0.5539757899927906
$ws.Range("C17").Value = 0.1909381264423828
$ws.Range("D17").Value = 0.05733679404445269
$ws.Range("F17").Value = 0.9670462914754907
$ws.Range("G17").Value = 0.002458675430808839
$ws.Range("I17").Value = 0.9342874829692285
$ws.Range("K17").Value = 0.3759038517319766
$ws.Range("L17").Value = 0.3172558355303039
$ws.Range("M17").Value = 0.1908652249031277
$ws.Range("O17").Value = 3.4181082293656
# Row 18
$ws.Range("B18").Value = 0.544017925263887
$ws.Range("C18").Value = 0.1906797022223401
$ws.Range("D18").Value = 0.05673633748222073
$ws.Range("F18").Value = 0.9668332013921557
$ws.Range("G18").Value = 0.002459211087827701
$ws.Range("I18").Value = 0.9357194303282128
$ws.Range("K18").Value = 0.3669164343047271
$ws.Range("L18").Value = 0.3161451864811653
$ws.Range("M18").Value = 0.1887327437310802
$ws.Range("O18").Value = 3.42110878756958
# Row 19
$ws.Range("B19").Value = 0.5406476316593682
$ws.Range("C19").Value = 0.1905922463300698
$ws.Range("D19").Value = 0.05653277302523918
$ws.Range("F19").Value = 0.9667700632640219
$ws.Range("G19").Value = 0.002459393740292446
$ws.Range("I19").Value = 0.9362117229186673
$ws.Range("K19").Value = 0.3638725734034267
$ws.Range("L19").Value = 0.3157717483243516
$ws.Range("M19").Value = 0.1880118142232519
$ws.Range("O19").Value = 3.422154478902655
# Row 20
$ws.Range("B20").Value = 0.5558193648984115
$ws.Range("C20").Value = 0.1909859747511504
$ws.Range("D20").Value = 0.05744780126929072
$ws.Range("F20").Value = 0.9670900093036892
$ws.Range("G20").Value = 0.002458576904171564
$ws.Range("I20").Value = 0.9340260049718268
$ws.Range("K20").Value = 0.3775668011252549
$ws.Range("L20").Value = 0.3174626302212999
$ws.Range("M20").Value = 0.1912604169878804
$ws.Range("O20").Value = 3.417567026161834
# Row 21
$ws.Range("B21").Value = 0.6068797222267222
$ws.Range("C21").Value = 0.1923116496538171
$ws.Range("D21").Value = 0.06050394214695132
$ws.Range("F21").Value = 0.9687915966313483
$ws.Range("G21").Value = 0.00245592183285077
$ws.Range("I21").Value = 0.9272011756170926
$ws.Range("K21").Value = 0.4235143466073623
$ws.Range("L21").Value = 0.3233248901143497
$ws.Range("M21").Value = 0.2022506502047534
$ws.Range("O21").Value = 3.404223462846943
# Row 22
$ws.Range("B22").Value = 0.6403109666751732
$ws.Range("C22").Value = 0.1931799929236107
$ws.Range("D22").Value = 0.06248756892084373
$ws.Range("F22").Value = 0.9703704645249047
$ws.Range("G22").Value = 0.002454252945676713
$ws.Range("I22").Value = 0.9231274189084573
$ws.Range("K22").Value = 0.4534939771348832
$ws.Range("L22").Value = 0.3272908201059153
$ws.Range("M22").Value = 0.2094888851678292
$ws.Range("O22").Value = 3.397045004103319
# Row 23
$ws.Range("B23").Value = 0.6224627726755045
$ws.Range("C23").Value = 0.1927163740238953
$ws.Range("D23").Value = 0.06143013675273323
$ws.Range("F23").Value = 0.9694850980626342
$ws.Range("G23").Value = 0.002455137607671701
$ws.Range("I23").Value = 0.9252662841134409
$ws.Range("K23").Value = 0.437498006187127
$ws.Range("L23").Value = 0.3251618376388166
$ws.Range("M23").Value = 0.2056206708975026
$ws.Range("O23").Value = 3.400735172808226
# Row 24
$ws.Range("B24").Value = 0.5549858760332995
$ws.Range("C24").Value = 0.1909643421446887
$ws.Range("D24").Value = 0.05739762048008146
$ws.Range("F24").Value = 0.9670700814520785
$ws.Range("G24").Value = 0.002458621423928997
$ws.Range("I24").Value = 0.9341440819292401
$ws.Range("K24").Value = 0.3768150104939991
$ws.Range("L24").Value = 0.3173690926565058
$ws.Range("M24").Value = 0.1910817339212869
$ws.Range("O24").Value = 3.41781116053491
# Row 25
$ws.Range("B25").Value = 0.4824814142274079
$ws.Range("C25").Value = 0.1890837799269534
$ws.Range("D25").Value = 0.05298988616410583
$ws.Range("F25").Value = 0.9664706666923308
$ws.Range("G25").Value = 0.002462665907407074
$ws.Range("I25").Value = 0.9453807581443883
$ws.Range("K25").Value = 0.3111619797895173
$ws.Range("L25").Value = 0.3095437646312575
$ws.Range("M25").Value = 0.1756417802006069
$ws.Range("O25").Value = 3.442853458381762
